# Apply the "Add data for 2023-12-05" update to the CTA violent-crime
# year-to-date workbook.
#
# The upstream commit adds one more day's worth of incident records to the
# rolling "year to date" counts. Because each neighborhood sheet already
# stores plain numeric totals (no formulas anywhere in the workbook), the
# edit shows up as a handful of cell values incrementing by 1-3, plus one
# brand-new cell (Fuller Park's Aggravated Battery / 2023 cell, which had no
# recorded incidents before and gets its first one). The Citywide Totals and
# By Neighborhood summary sheets roll the same deltas up to their own
# Robbery/Total rows.

$wb = $excel.ActiveWorkbook

function Set-Cell {
    param(
        [string]$SheetName,
        [string]$CellRef,
        [double]$Value
    )
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Range($CellRef).Value = $Value
}

# --- Citywide Totals ---------------------------------------------------
Set-Cell "Citywide Totals" "J2" 121
Set-Cell "Citywide Totals" "E3" 144
Set-Cell "Citywide Totals" "J3" 231
Set-Cell "Citywide Totals" "B6" 374
Set-Cell "Citywide Totals" "C6" 475
Set-Cell "Citywide Totals" "D6" 413
Set-Cell "Citywide Totals" "E6" 470
Set-Cell "Citywide Totals" "F6" 531
Set-Cell "Citywide Totals" "I6" 498
Set-Cell "Citywide Totals" "J6" 415
Set-Cell "Citywide Totals" "B7" 499
Set-Cell "Citywide Totals" "C7" 630
Set-Cell "Citywide Totals" "D7" 644
Set-Cell "Citywide Totals" "E7" 695
Set-Cell "Citywide Totals" "F7" 767
Set-Cell "Citywide Totals" "I7" 831
Set-Cell "Citywide Totals" "J7" 788

# --- By Neighborhood (per-neighborhood Total rows + grand Total) -------
Set-Cell "By Neighborhood" "F7" 11
Set-Cell "By Neighborhood" "B8" 31
Set-Cell "By Neighborhood" "I21" 15
Set-Cell "By Neighborhood" "J29" 13
Set-Cell "By Neighborhood" "E32" 66
Set-Cell "By Neighborhood" "J32" 46
Set-Cell "By Neighborhood" "C36" 39
Set-Cell "By Neighborhood" "E53" 82
Set-Cell "By Neighborhood" "F53" 82
Set-Cell "By Neighborhood" "D65" 25
Set-Cell "By Neighborhood" "J78" 6
Set-Cell "By Neighborhood" "E91" 7
Set-Cell "By Neighborhood" "B96" 16
Set-Cell "By Neighborhood" "B98" 499
Set-Cell "By Neighborhood" "C98" 630
Set-Cell "By Neighborhood" "D98" 644
Set-Cell "By Neighborhood" "E98" 695
Set-Cell "By Neighborhood" "F98" 767
Set-Cell "By Neighborhood" "I98" 831
Set-Cell "By Neighborhood" "J98" 788

# --- Individual neighborhood sheets ------------------------------------
Set-Cell "Auburn Gresham" "F5" 7
Set-Cell "Auburn Gresham" "F6" 11

Set-Cell "Austin" "B5" 22
Set-Cell "Austin" "B6" 31

Set-Cell "Chinatown" "I6" 9
Set-Cell "Chinatown" "I7" 15

Set-Cell "Garfield Park" "J3" 14
Set-Cell "Garfield Park" "E6" 53
Set-Cell "Garfield Park" "J6" 30
Set-Cell "Garfield Park" "E7" 66
Set-Cell "Garfield Park" "J7" 46

Set-Cell "Grand Crossing" "C6" 34
Set-Cell "Grand Crossing" "C7" 39

Set-Cell "Rush & Division" "J4" 5
Set-Cell "Rush & Division" "J5" 6

Set-Cell "Woodlawn" "B5" 11
Set-Cell "Woodlawn" "B6" 16

Set-Cell "Loop" "E3" 13
Set-Cell "Loop" "F6" 61
Set-Cell "Loop" "E7" 82
Set-Cell "Loop" "F7" 82

Set-Cell "West Loop" "E6" 6
Set-Cell "West Loop" "E7" 7

Set-Cell "North Lawndale" "D5" 24
Set-Cell "North Lawndale" "D6" 25

# Fuller Park gains its first-ever 2023 Aggravated Battery incident (new
# cell J2), which also bumps the sheet's own Total row (J6).
Set-Cell "Fuller Park" "J2" 1
Set-Cell "Fuller Park" "J6" 13
